$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..49: [date, B, C, D]
# (diff reorders each 12-row year block so Oct/Nov/Dec precede Jan-Sep)
$data = @(
  @("2014-10", 100.037, 102.4351, 102.4479),
  @("2014-11", 100, 102.6243, 101.9211),
  @("2014-12", 100, 102.8596, 102.1653),
  @("2014-01", 100.037, 100.0364, 102.6889),
  @("2014-02", 100.037, 100.1379, 102.8531),
  @("2014-03", 100.037, 99.8967, 102.5829),
  @("2014-04", 100.037, 100.0517, 102.6985),
  @("2014-05", 100.037, 100.8555, 103.019),
  @("2014-06", 100.037, 101.0447, 103.0387),
  @("2014-07", 100.037, 101.0579, 103.0811),
  @("2014-08", 100.037, 101.7929, 103.1783),
  @("2014-09", 100.037, 102.4299, 103.0798),
  @("2015-10", 95, 99.40000000000001, 102.3),
  @("2015-11", 95, 98.8458, 102.2828),
  @("2015-12", 95, 98.6323, 102.0785),
  @("2015-01", 100, 103.093, 102.4015),
  @("2015-02", 100, 102.6343, 102.1907),
  @("2015-03", 100, 102.1288, 102.357),
  @("2015-04", 100, 102.1288, 102.4084),
  @("2015-05", 100, 101.3743, 102.214),
  @("2015-06", 100, 101.2307, 102.5949),
  @("2015-07", 100, 101.1896, 102.5619),
  @("2015-08", 100, 100.0674, 102.4717),
  @("2015-09", 100, 99.30929999999999, 102.3147),
  @("2016-10", 99.40000000000001, 99.5, 102.7),
  @("2016-11", 99.3, 99.59999999999999, 102.6),
  @("2016-12", 100, 99.59999999999999, 102.7),
  @("2016-01", 100.2999, 99.4482, 102.1153),
  @("2016-02", 100.2697, 99.60290000000001, 102.3254),
  @("2016-03", 99.827, 99.4572, 102.1788),
  @("2016-04", 100.058, 99.53279999999999, 102.6876),
  @("2016-05", 99.90000000000001, 99, 103.1),
  @("2016-06", 99.8, 99.09999999999999, 102.9),
  @("2016-07", 99.5, 98.59999999999999, 102.9),
  @("2016-08", 99.59999999999999, 98.7, 102.7),
  @("2016-09", 99.59999999999999, 99, 102.7),
  @("2017-10", 100, 102.6, 102),
  @("2017-11", 100, 102.5, 102.1),
  @("2017-12", 100, 101.9, 101.9),
  @("2017-01", 100, 100.3, 102.3),
  @("2017-02", 100, 100.7, 102.3),
  @("2017-03", 100, 101.8, 102.4),
  @("2017-04", 100, 102.1, 101.8),
  @("2017-05", 100, 102.9, 101.4),
  @("2017-06", 100, 102.7, 101.5),
  @("2017-07", 100, 103.3, 101.4),
  @("2017-08", 100, 103.4, 101.6),
  @("2017-09", 100, 103.4, 101.9)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
